$d = $word.ActiveDocument

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Contains("视觉显著性，作为一种")) {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Target paragraph not found"
}

$start = $target.Range.Start
$end = $target.Range.End - 1
$r = $d.Range($start, $end)

$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:color w:val="5B9BD5" w:themeColor="accent1"/></w:rPr><w:t>视觉显著性</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:color w:val="5B9BD5" w:themeColor="accent1"/></w:rPr><w:t>，作为一种</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:color w:val="5B9BD5" w:themeColor="accent1"/></w:rPr><w:t>重要的视觉特征，最近几年已经得到越来越多的关注。</w:t></w:r><w:r w:rsidR="00394A1B"><w:rPr><w:color w:val="5B9BD5" w:themeColor="accent1"/></w:rPr><w:t>[1-4, 26, 27]</w:t></w:r><w:r w:rsidR="000A3027"><w:rPr><w:rFonts w:hint="eastAsia"/><w:color w:val="5B9BD5" w:themeColor="accent1"/></w:rPr><w:t>视觉显著性表明</w:t></w:r><w:bookmarkStart w:id="0" w:name="OLE_LINK1"/><w:bookmarkStart w:id="1" w:name="OLE_LINK2"/><w:r w:rsidR="000A3027"><w:rPr><w:rFonts w:hint="eastAsia"/><w:color w:val="5B9BD5" w:themeColor="accent1"/></w:rPr><w:t>子区域是否与周围环境有明显差异兵器能快速</w:t></w:r><w:r w:rsidR="00FE7F22"><w:rPr><w:rFonts w:hint="eastAsia"/><w:color w:val="5B9BD5" w:themeColor="accent1"/></w:rPr><w:t>引起</w:t></w:r><w:r w:rsidR="000A3027"><w:rPr><w:rFonts w:hint="eastAsia"/><w:color w:val="5B9BD5" w:themeColor="accent1"/></w:rPr><w:t>观察用户的注意</w:t></w:r><w:bookmarkStart w:id="2" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:bookmarkEnd w:id="1"/><w:bookmarkEnd w:id="2"/><w:r w:rsidR="00FE7F22"><w:rPr><w:rFonts w:hint="eastAsia"/><w:color w:val="5B9BD5" w:themeColor="accent1"/></w:rPr><w:t>。自从</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:color w:val="5B9BD5" w:themeColor="accent1"/></w:rPr><w:t>Itti</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00FE7F22"><w:rPr><w:color w:val="5B9BD5" w:themeColor="accent1"/></w:rPr><w:t>[7]</w:t></w:r><w:r w:rsidR="00FE7F22"><w:rPr><w:rFonts w:hint="eastAsia"/><w:color w:val="5B9BD5" w:themeColor="accent1"/></w:rPr><w:t>提出并应用</w:t></w:r><w:r w:rsidR="00674278"><w:rPr><w:rFonts w:hint="eastAsia"/><w:color w:val="5B9BD5" w:themeColor="accent1"/></w:rPr><w:t>计算机视觉的视觉注意力理论以来</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:color w:val="5B9BD5" w:themeColor="accent1"/></w:rPr><w:t>，所有</w:t></w:r><w:r w:rsidR="00674278"><w:rPr><w:rFonts w:hint="eastAsia"/><w:color w:val="5B9BD5" w:themeColor="accent1"/></w:rPr><w:t>后来</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:color w:val="5B9BD5" w:themeColor="accent1"/></w:rPr><w:t>提出的方法都是</w:t></w:r><w:r w:rsidR="00674278"><w:rPr><w:rFonts w:hint="eastAsia"/><w:color w:val="5B9BD5" w:themeColor="accent1"/></w:rPr><w:t>遵循</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:color w:val="5B9BD5" w:themeColor="accent1"/></w:rPr><w:t>文献[3]中</w:t></w:r><w:r w:rsidR="00674278"><w:rPr><w:rFonts w:hint="eastAsia"/><w:color w:val="5B9BD5" w:themeColor="accent1"/></w:rPr><w:t>总结</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:color w:val="5B9BD5" w:themeColor="accent1"/></w:rPr><w:t>的人类视觉注意</w:t></w:r><w:r w:rsidR="00674278"><w:rPr><w:rFonts w:hint="eastAsia"/><w:color w:val="5B9BD5" w:themeColor="accent1"/></w:rPr><w:t>四个基本原则</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:color w:val="5B9BD5" w:themeColor="accent1"/></w:rPr><w:t>中的</w:t></w:r><w:r w:rsidR="00674278"><w:rPr><w:rFonts w:hint="eastAsia"/><w:color w:val="5B9BD5" w:themeColor="accent1"/></w:rPr><w:t>一个或者几个，这些原则考虑了</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:color w:val="5B9BD5" w:themeColor="accent1"/></w:rPr><w:t>局部底层</w:t></w:r><w:r w:rsidR="00674278"><w:rPr><w:rFonts w:hint="eastAsia"/><w:color w:val="5B9BD5" w:themeColor="accent1"/></w:rPr><w:t>因素，全局因素，视觉组织规则以及高级因素</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:color w:val="5B9BD5" w:themeColor="accent1"/></w:rPr><w:t>。</w:t></w:r><w:r w:rsidR="00674278"><w:rPr><w:rFonts w:hint="eastAsia"/><w:color w:val="5B9BD5" w:themeColor="accent1"/></w:rPr><w:t>视觉显著性可以被用到许多计算机视觉任务</w:t></w:r><w:r w:rsidR="00674278"><w:rPr><w:color w:val="5B9BD5" w:themeColor="accent1"/></w:rPr><w:t>[</w:t></w:r><w:r w:rsidR="00674278"><w:rPr><w:rFonts w:hint="eastAsia"/><w:color w:val="5B9BD5" w:themeColor="accent1"/></w:rPr><w:t>28</w:t></w:r><w:r w:rsidR="00674278"><w:rPr><w:color w:val="5B9BD5" w:themeColor="accent1"/></w:rPr><w:t>]</w:t></w:r><w:r w:rsidR="00674278"><w:rPr><w:rFonts w:hint="eastAsia"/><w:color w:val="5B9BD5" w:themeColor="accent1"/></w:rPr><w:t>或者图像处理</w:t></w:r><w:r w:rsidR="00674278"><w:rPr><w:color w:val="5B9BD5" w:themeColor="accent1"/></w:rPr><w:t>[29]</w:t></w:r><w:r w:rsidR="00674278"><w:rPr><w:rFonts w:hint="eastAsia"/><w:color w:val="5B9BD5" w:themeColor="accent1"/></w:rPr><w:t>中。</w:t></w:r><w:r w:rsidR="008F45EB"><w:rPr><w:rFonts w:hint="eastAsia"/><w:color w:val="5B9BD5" w:themeColor="accent1"/></w:rPr><w:t>许多机器人系统也利用视觉显著性去进行对象识别与检测</w:t></w:r><w:r w:rsidR="008F45EB"><w:rPr><w:color w:val="5B9BD5" w:themeColor="accent1"/></w:rPr><w:t>[5, 6]</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:color w:val="5B9BD5" w:themeColor="accent1"/></w:rPr><w:t>。</w:t></w:r></w:p>
'@

$r.InsertXML($xml) | Out-Null
Write-Host "DONE"
